# Update 2024 year-to-date (column K) violent crime counts for 2024-07-10
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 4096
$ws.Range("K3").Value = 4207
$ws.Range("I4").Value = 1796
$ws.Range("J4").Value = 1822
$ws.Range("K4").Value = 847
$ws.Range("K5").Value = 300
$ws.Range("K6").Value = 4704
$ws.Range("I7").Value = 26251
$ws.Range("J7").Value = 29291
$ws.Range("K7").Value = 14154

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 46
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 270
$ws.Range("K3").Value = 288
$ws.Range("K6").Value = 322
$ws.Range("K7").Value = 961

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 300

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 159
$ws.Range("K3").Value = 222
$ws.Range("K6").Value = 169
$ws.Range("K7").Value = 588

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 79
$ws.Range("K3").Value = 86
$ws.Range("K7").Value = 244

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 411
$ws.Range("K8").Value = 961
$ws.Range("K9").Value = 59
$ws.Range("K11").Value = 278
$ws.Range("K14").Value = 75
$ws.Range("K18").Value = 97
$ws.Range("K19").Value = 436
$ws.Range("K20").Value = 317
$ws.Range("K21").Value = 41
$ws.Range("K23").Value = 142
$ws.Range("K27").Value = 138
$ws.Range("K29").Value = 741
$ws.Range("K33").Value = 588
$ws.Range("K36").Value = 178
$ws.Range("K41").Value = 118
$ws.Range("K42").Value = 513
$ws.Range("K43").Value = 126
$ws.Range("K44").Value = 129
$ws.Range("K46").Value = 32
$ws.Range("K49").Value = 83
$ws.Range("K51").Value = 174
$ws.Range("K52").Value = 387
$ws.Range("K53").Value = 191
$ws.Range("K54").Value = 263
$ws.Range("I63").Value = 213
$ws.Range("J63").Value = 106
$ws.Range("K63").Value = 46
$ws.Range("K64").Value = 86
$ws.Range("K67").Value = 550
$ws.Range("K70").Value = 24
$ws.Range("K72").Value = 66
$ws.Range("K75").Value = 47
$ws.Range("K76").Value = 201
$ws.Range("K78").Value = 168
$ws.Range("K79").Value = 365
$ws.Range("K83").Value = 300
$ws.Range("K85").Value = 634
$ws.Range("K86").Value = 96
$ws.Range("K88").Value = 162
$ws.Range("K89").Value = 199
$ws.Range("K91").Value = 157
$ws.Range("K93").Value = 51
$ws.Range("K95").Value = 244
$ws.Range("K97").Value = 121
$ws.Range("I101").Value = 26251
$ws.Range("J101").Value = 29291
$ws.Range("K101").Value = 14154

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 161
$ws.Range("K5").Value = 10
$ws.Range("K7").Value = 550

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 129
$ws.Range("K7").Value = 263

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K6").Value = 207
$ws.Range("K7").Value = 741

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 134
$ws.Range("K7").Value = 436

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 139
$ws.Range("K3").Value = 165
$ws.Range("K4").Value = 22
$ws.Range("K6").Value = 184
$ws.Range("K7").Value = 513

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 73
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K4").Value = 24
$ws.Range("K7").Value = 365

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 98
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 317

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 29
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 70
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 147
$ws.Range("K3").Value = 135
$ws.Range("K5").Value = 17
$ws.Range("K7").Value = 411

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 90
$ws.Range("K4").Value = 15
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K2").Value = 9
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 53
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 46
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 174

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 212
$ws.Range("K5").Value = 18
$ws.Range("K6").Value = 148
$ws.Range("K7").Value = 634

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 101
$ws.Range("K7").Value = 387

Write-Output "Updated 159 cells across 43 worksheets."